$d = $word.ActiveDocument

# --- Change 2 (later in the document): remove the empty paragraph that
# follows "Prioridade: Baixa" (it gets merged into the preceding
# paragraph, which keeps only the trailing empty run that used to belong
# to the now-removed paragraph). We apply this edit first since it is
# further down the document, so it does not shift the character offsets
# used for the first change below.
$rngBaixa = $d.Content
$rngBaixa.Find.Execute("Prioridade: Baixa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterBaixa = $rngBaixa.End
$d.Range($afterBaixa + 1, $afterBaixa + 2).Delete()

# --- Change 1 (earlier in the document): remove the blank, left-aligned
# paragraph that immediately precedes the "Projeto: ..." title paragraph,
# and extend the title text.
$rngTitle = $d.Content
$rngTitle.Find.Execute("Projeto: Sistema WEB para gerenciamento", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleStart = $rngTitle.Start
$d.Range($titleStart - 1, $titleStart).Delete()

$d.Content.Find.Execute("Projeto: Sistema WEB para gerenciamento", $true, $false, $false, $false, $false, $true, 1, $false, "Projeto: Sistema WEB para gerenciamento de assistência técnica a computadores", 2)
